$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 82.71429000000001
$ws.Range("I8").Value = 25.8
$ws.Range("K8").Value = 77.40000000000001
$ws.Range("M8").Value = 61.59999999999999
$ws.Range("H17").Value = 1005.403
$ws.Range("J17").Value = 1005.1754
$ws.Range("L17").Value = 3015.5262
$ws.Range("N17").Value = -3351.5262
$ws.Range("H64").Value = 3335.5
$ws.Range("I64").Value = 3000
$ws.Range("K64").Value = 3000
$ws.Range("M64").Value = -2752
$ws.Range("H67").Value = 3335.5
$ws.Range("I67").Value = 3000
$ws.Range("K67").Value = 3000
$ws.Range("M67").Value = -2142
$ws.Range("H96").Value = 2056.0588
$ws.Range("I96").Value = 660.6667
$ws.Range("K96").Value = 1982.0001
$ws.Range("M96").Value = -609.0001
$ws.Range("H137").Value = 3870.6272
$ws.Range("I137").Value = 2001.5686
$ws.Range("K137").Value = 6004.7058
$ws.Range("M137").Value = -3454.7058
$ws.Range("H138").Value = 2565.8372
$ws.Range("I138").Value = 1944.6571
$ws.Range("J138").Value = 5283.5
$ws.Range("K138").Value = 5833.971299999999
$ws.Range("L138").Value = 15850.5
$ws.Range("M138").Value = -693.9712999999992
$ws.Range("N138").Value = -26130.5
$ws.Range("H141").Value = 5130.2974
$ws.Range("I141").Value = 2176.1072
$ws.Range("K141").Value = 6528.321599999999
$ws.Range("M141").Value = -1348.321599999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3104.5454
$ws.Range("I86").Value = 3100.5
$ws.Range("K86").Value = 3100.5
$ws.Range("M86").Value = -1977.5
$ws.Range("H89").Value = 3104.5454
$ws.Range("I89").Value = 3100.5
$ws.Range("K89").Value = 15502.5
$ws.Range("M89").Value = -9886.5
$ws.Range("H107").Value = 1337.5333
$ws.Range("I107").Value = 1051.1818
$ws.Range("K107").Value = 1051.1818
$ws.Range("M107").Value = 868.8181999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4185.615
$ws.Range("I31").Value = 3041.3
$ws.Range("K31").Value = 3041.3
$ws.Range("M31").Value = -2746.3
$ws.Range("H34").Value = 4185.615
$ws.Range("I34").Value = 3041.3
$ws.Range("K34").Value = 3041.3
$ws.Range("M34").Value = -2839.3
$ws.Range("H80").Value = 21719
$ws.Range("J80").Value = 21719
$ws.Range("L80").Value = 21719
$ws.Range("N80").Value = -23965
$ws.Range("H83").Value = 21719
$ws.Range("J83").Value = 21719
$ws.Range("L83").Value = 65157
$ws.Range("N83").Value = -76389

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 56.7
$ws.Range("I40").Value = 58
$ws.Range("J40").Value = 55.4
$ws.Range("K40").Value = 232
$ws.Range("L40").Value = 221.6
$ws.Range("M40").Value = -163
$ws.Range("N40").Value = -359.6
$ws.Range("H82").Value = 14802.6
$ws.Range("H85").Value = 14802.6
$ws.Range("H98").Value = 833.375
$ws.Range("I98").Value = 1223.25
$ws.Range("J98").Value = 443.5
$ws.Range("K98").Value = 3669.75
$ws.Range("L98").Value = 1330.5
$ws.Range("M98").Value = -2171.75
$ws.Range("N98").Value = -4326.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 16687.6
$ws.Range("J15").Value = 16687.6
$ws.Range("L15").Value = 16687.6
$ws.Range("N15").Value = -17263.6
$ws.Range("H42").Value = 60054
$ws.Range("J42").Value = 59949.5
$ws.Range("L42").Value = 59949.5
$ws.Range("N42").Value = -60919.5
$ws.Range("H70").Value = 21112.23
$ws.Range("I70").Value = 28682.75
$ws.Range("J70").Value = 8999.4
$ws.Range("K70").Value = 28682.75
$ws.Range("L70").Value = 8999.4
$ws.Range("M70").Value = -28412.75
$ws.Range("N70").Value = -9539.4
$ws.Range("H73").Value = 21112.23
$ws.Range("I73").Value = 28682.75
$ws.Range("J73").Value = 8999.4
$ws.Range("K73").Value = 28682.75
$ws.Range("L73").Value = 8999.4
$ws.Range("M73").Value = -27746.75
$ws.Range("N73").Value = -10871.4
$ws.Range("H81").Value = 16687.6
$ws.Range("J81").Value = 16687.6
$ws.Range("L81").Value = 16687.6
$ws.Range("N81").Value = -18683.6
$ws.Range("H84").Value = 16687.6
$ws.Range("J84").Value = 16687.6
$ws.Range("L84").Value = 50062.8
$ws.Range("N84").Value = -60046.8
$ws.Range("H97").Value = 1067.4
$ws.Range("I97").Value = 1150.9166
$ws.Range("K97").Value = 1150.9166
$ws.Range("M97").Value = -654.9166
$ws.Range("H98").Value = 34428.668
$ws.Range("J98").Value = 34428.668
$ws.Range("L98").Value = 34428.668
$ws.Range("N98").Value = -40418.668
$ws.Range("H99").Value = 18293
$ws.Range("I99").Value = 5821.6665
$ws.Range("J99").Value = 37000
$ws.Range("K99").Value = 5821.6665
$ws.Range("L99").Value = 37000
$ws.Range("M99").Value = -3575.6665
$ws.Range("N99").Value = -41492
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").Value = $null
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").Value = $null
$ws.Range("H107").Value = 296.10526
$ws.Range("I107").Value = 346.07693
$ws.Range("K107").Value = 346.07693
$ws.Range("M107").Value = 1573.92307
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").Value = $null
$ws.Range("H113").Value = 1769.7
$ws.Range("I113").Value = 1758.5294
$ws.Range("K113").Value = 1758.5294
$ws.Range("M113").Value = 411.4706000000001
$ws.Range("H115").Value = 60054
$ws.Range("J115").Value = 59949.5
$ws.Range("L115").Value = 59949.5
$ws.Range("N115").Value = -62299.5
$ws.Range("H117").Value = 91452.75
$ws.Range("J117").Value = 91452.75
$ws.Range("L117").Value = 91452.75
$ws.Range("N117").Value = -98336.75
$ws.Range("H118").Value = 17732.834
$ws.Range("J118").Value = 17732.834
$ws.Range("L118").Value = 17732.834
$ws.Range("N118").Value = -21046.834
$ws.Range("H120").Value = 38210
$ws.Range("J120").Value = 38210
$ws.Range("L120").Value = 38210
$ws.Range("N120").Value = -47886
$ws.Range("H121").Value = 40500
$ws.Range("J121").Value = 40500
$ws.Range("L121").Value = 40500
$ws.Range("N121").Value = -43994
$ws.Range("H122").Value = 2194.348
$ws.Range("I122").Value = 1706.6428
$ws.Range("K122").Value = 5119.928400000001
$ws.Range("M122").Value = -2669.928400000001
$ws.Range("H123").Value = 50585.2
$ws.Range("J123").Value = 50585.2
$ws.Range("L123").Value = 50585.2
$ws.Range("N123").Value = -55485.2
$ws.Range("H126").Value = 2529
$ws.Range("I126").Value = 2436.2
$ws.Range("K126").Value = 7308.599999999999
$ws.Range("M126").Value = -4838.599999999999
$ws.Range("H127").Value = 93425.336
$ws.Range("J127").Value = 93425.336
$ws.Range("L127").Value = 93425.336
$ws.Range("N127").Value = -103345.336
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").Value = $null
$ws.Range("H129").Value = 59000
$ws.Range("J129").Value = 59000
$ws.Range("L129").Value = 59000
$ws.Range("N129").Value = -69000
$ws.Range("H130").Value = 49325
$ws.Range("J130").Value = 49325
$ws.Range("L130").Value = 49325
$ws.Range("N130").Value = -59365
$ws.Range("H132").Value = 20348.125
$ws.Range("I132").Value = 21326.428
$ws.Range("K132").Value = 63979.284
$ws.Range("M132").Value = -61449.284

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 11892.211
$ws.Range("I68").Value = 10380.385
$ws.Range("K68").Value = 10380.385
$ws.Range("M68").Value = -9631.385
$ws.Range("H71").Value = 11892.211
$ws.Range("I71").Value = 10380.385
$ws.Range("K71").Value = 51901.925
$ws.Range("M71").Value = -48157.925
$ws.Range("H80").Value = 45000
$ws.Range("J80").Value = 45000
$ws.Range("L80").Value = 45000
$ws.Range("N80").Value = -47246
$ws.Range("H83").Value = 45000
$ws.Range("J83").Value = 45000
$ws.Range("L83").Value = 135000
$ws.Range("N83").Value = -146232
$ws.Range("H136").Value = 1842.3334
$ws.Range("I136").Value = 1589.72
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 4769.16
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -2219.16
$ws.Range("N136").Value = -20100

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 147421.58
$ws.Range("J81").Value = 204290
$ws.Range("L81").Value = 408580
$ws.Range("N81").Value = -410702
$ws.Range("H84").Value = 147421.58
$ws.Range("J84").Value = 204290
$ws.Range("L84").Value = 2042900
$ws.Range("N84").Value = -2053508
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = $null
$ws.Range("H107").Value = 2307.5588
$ws.Range("I107").Value = 1086.2916
$ws.Range("J107").Value = 5238.6
$ws.Range("K107").Value = 3258.8748
$ws.Range("L107").Value = 15715.8
$ws.Range("M107").Value = -1338.8748
$ws.Range("N107").Value = -19555.8
$ws.Range("H132").Value = 2914.4583
$ws.Range("I132").Value = 1977.6428
$ws.Range("K132").Value = 5932.928400000001
$ws.Range("M132").Value = -3402.928400000001
$ws.Range("H136").Value = 1517.6086
$ws.Range("I136").Value = 1495.4762
$ws.Range("K136").Value = 4486.4286
$ws.Range("M136").Value = -1936.4286
